$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 10-13 (old rows for "Resolving-Mac" A-column group no longer present / table shrinks to 9 rows)
$ws.Rows("10:13").Delete() | Out-Null

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("D2").Value = "MuSCs"
$ws.Range("G2").Value = 169.915657
$ws.Range("H2").Value = 509.746971
$ws.Range("I2").Value = 0.4441184931734509
$ws.Range("J2").Value = 0.4441184931734509
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.5092303333333333
$ws.Range("N2").Value = 1.527691
$ws.Range("O2").Value = 0.02338915669123285
$ws.Range("P2").Value = 0.02338915669123285
$ws.Range("Q2").Value = 86.52620665266234
$ws.Range("R2").Value = 778.7358598739611
$ws.Range("S2").Value = 0.01038755702630807
$ws.Range("T2").Value = 0.01038755702630807

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("G3").Value = 169.915657
$ws.Range("H3").Value = 509.746971
$ws.Range("I3").Value = 0.4441184931734509
$ws.Range("J3").Value = 0.4441184931734509
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 21.26283866666667
$ws.Range("N3").Value = 63.788516
$ws.Range("O3").Value = 0.9766108433087671
$ws.Range("P3").Value = 0.9766108433087671
$ws.Range("Q3").Value = 3612.889201731671
$ws.Range("R3").Value = 32516.00281558504
$ws.Range("S3").Value = 0.4337309361471428
$ws.Range("T3").Value = 0.4337309361471429

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("G4").Value = 68.382243
$ws.Range("H4").Value = 205.146729
$ws.Range("I4").Value = 0.1787346690539575
$ws.Range("J4").Value = 0.1787346690539575
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.5092303333333333
$ws.Range("N4").Value = 1.527691
$ws.Range("O4").Value = 0.02338915669123285
$ws.Range("P4").Value = 0.02338915669123285
$ws.Range("Q4").Value = 34.822312396971
$ws.Range("R4").Value = 313.400811572739
$ws.Range("S4").Value = 0.00418045318065866
$ws.Range("T4").Value = 0.00418045318065866

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("G5").Value = 68.382243
$ws.Range("H5").Value = 205.146729
$ws.Range("I5").Value = 0.1787346690539575
$ws.Range("J5").Value = 0.1787346690539575
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 21.26283866666667
$ws.Range("N5").Value = 63.788516
$ws.Range("O5").Value = 0.9766108433087671
$ws.Range("P5").Value = 0.9766108433087671
$ws.Range("Q5").Value = 1454.000600573796
$ws.Range("R5").Value = 13086.00540516416
$ws.Range("S5").Value = 0.1745542158732989
$ws.Range("T5").Value = 0.1745542158732989

# Row 6
$ws.Range("A6").Value = "MuSCs"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("G6").Value = 53.27463399999999
$ws.Range("H6").Value = 159.823902
$ws.Range("I6").Value = 0.1392470275793777
$ws.Range("J6").Value = 0.1392470275793778
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.5092303333333333
$ws.Range("N6").Value = 1.527691
$ws.Range("O6").Value = 0.02338915669123285
$ws.Range("P6").Value = 0.02338915669123285
$ws.Range("Q6").Value = 27.12905963003133
$ws.Range("R6").Value = 244.161536670282
$ws.Range("S6").Value = 0.003256870546842488
$ws.Range("T6").Value = 0.003256870546842489

# Row 7
$ws.Range("A7").Value = "MuSCs"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("G7").Value = 53.27463399999999
$ws.Range("H7").Value = 159.823902
$ws.Range("I7").Value = 0.1392470275793777
$ws.Range("J7").Value = 0.1392470275793778
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 21.26283866666667
$ws.Range("N7").Value = 63.788516
$ws.Range("O7").Value = 0.9766108433087671
$ws.Range("P7").Value = 0.9766108433087671
$ws.Range("Q7").Value = 1132.769947767714
$ws.Range("R7").Value = 10194.92952990943
$ws.Range("S7").Value = 0.1359901570325353
$ws.Range("T7").Value = 0.1359901570325353

# Row 8
$ws.Range("A8").Value = "Resolving-Mac"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("G8").Value = 91.01828266666666
$ws.Range("H8").Value = 273.054848
$ws.Range("I8").Value = 0.2378998101932138
$ws.Range("J8").Value = 0.2378998101932138
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.5092303333333333
$ws.Range("N8").Value = 1.527691
$ws.Range("O8").Value = 0.02338915669123285
$ws.Range("P8").Value = 0.02338915669123285
$ws.Range("Q8").Value = 46.34927042177422
$ws.Range("R8").Value = 417.143433795968
$ws.Range("S8").Value = 0.005564275937423631
$ws.Range("T8").Value = 0.005564275937423632

# Row 9
$ws.Range("A9").Value = "Resolving-Mac"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("G9").Value = 91.01828266666666
$ws.Range("H9").Value = 273.054848
$ws.Range("I9").Value = 0.2378998101932138
$ws.Range("J9").Value = 0.2378998101932138
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 21.26283866666667
$ws.Range("N9").Value = 63.788516
$ws.Range("O9").Value = 0.9766108433087671
$ws.Range("P9").Value = 0.9766108433087671
$ws.Range("Q9").Value = 1935.307060058396
$ws.Range("R9").Value = 17417.76354052557
$ws.Range("S9").Value = 0.2323355342557901
$ws.Range("T9").Value = 0.2323355342557902
